# Generate Report for handback
#
# Marks the zh-cn / de-de handoff rows as "handed back": updates the
# Status text, records the handback timestamp, and fills in the
# "Latest Target File" / "Latest Handback File" columns (with their
# hyperlinks) for the single data row in each language sheet.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$mdFileName  = "d1fdf0e2-4eed-4a09-8182-aaca28568f5b.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/3fa0d1adafa078f899c39c1ac049c06c225a5a1d/e2e/d1fdf0e2-4eed-4a09-8182-aaca28568f5b.md"

$zhXlfFileName = "d1fdf0e2-4eed-4a09-8182-aaca28568f5b.2392d17c40969149aa70f06d1b0a302cdb6cd1fb.zh-cn.xlf"
$zhXlfUrl      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b5c62492f521b9c392dceb64493de3bc1047573/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/d1fdf0e2-4eed-4a09-8182-aaca28568f5b.2392d17c40969149aa70f06d1b0a302cdb6cd1fb.zh-cn.xlf"

$deXlfFileName = "d1fdf0e2-4eed-4a09-8182-aaca28568f5b.2392d17c40969149aa70f06d1b0a302cdb6cd1fb.de-de.xlf"
$deXlfUrl      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66253aef1e7ba44b0a17cf00f765312913cc253f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/d1fdf0e2-4eed-4a09-8182-aaca28568f5b.2392d17c40969149aa70f06d1b0a302cdb6cd1fb.de-de.xlf"

$zhHandbackDatetime = "2016-01-13 15:39:35"
$deHandbackDatetime = "2016-01-13 15:39:49"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusHandedBack

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, "", "", $mdFileName)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl, "", "", $zhXlfFileName)

$wsZh.Range("G2").Value = $zhHandbackDatetime

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusHandedBack

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, "", "", $mdFileName)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl, "", "", $deXlfFileName)

$wsDe.Range("G2").Value = $deHandbackDatetime
